# Updates the "cryptos" worksheet with refreshed price/volume data.
# Generated from the OOXML diff describing the crypto price update commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.449.71'
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = '2.254.36'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("E4").Value = '  +0.05%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '246.38'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E6").Value = '  +0.52%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '77.04'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("E8").Value = '  +0.02%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.625'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '45.28'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +10.75%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0953'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.06%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '7.29'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("E13").Value = '  -0.71%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '14.68'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.67%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.862'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").Value = '2.272.71'
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("D17").Value = '42.354.02'
$ws.Range("E18").Value = '  +4.27%  '
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("E21").Value = '  +3.11%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '232.19'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.92%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '9.14'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +26.71%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '11.58'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +4.55%  '
$ws.Range("E26").Value = '  -2.30%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '2.31'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  +1.88%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '167.58'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.12%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '20.71'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.34%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.0828'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.65%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '31.31'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -4.96%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.119'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '5.35'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +11.11%  '
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("E36").Value = '  -0.48%  '
$ws.Range("E37").Value = '  +6.28%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '14.37'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +6.50%  '
$ws.Range("E39").Value = '  +0.52%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '5.82'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.39%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '64.05'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +7.08%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.203'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.34%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '108.31'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -2.77%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '8.84'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("E45").Value = '  +3.30%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E47").Value = '  +4.92%  '
$ws.Range("E48").Value = '  +0.69%  '
$ws.Range("E49").Value = '  +2.26%  '
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("E51").Value = '  +0.85%  '
